$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- Update row 2 (data row for the new "inversion2" test case) ---
$ws.Range("D2").Value = "inversion2"
$ws.Range("E2").Value = "1234"
$ws.Range("O2").Value = "Prestamo personal ta"
$ws.Range("P2").Value = "29281005510"
$ws.Range("Q2").Value = "Pago total"
$ws.Range("R2").Value = "20561111"
$ws.Range("T2").Value = "406-130790-01"

# Re-assert the original cell formatting (quote-prefixed text style) for the
# two cells whose value-only assignment above could otherwise drop it.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null

# --- Row 3 keeps its data but swaps tipoPago with row 2 ("Pago total" <-> "Pago minimo") ---
$ws.Range("Q3").Value = "Pago mínimo"

# The tipoPago column carries a left-border accent that follows the "Pago minimo" row;
# move it from row 2 to row 3 to match the swapped values.
$ws.Range("Q2").Borders.Item(7).LineStyle = -4142
$ws.Range("Q3").Borders.Item(7).LineStyle = 1

# --- Update the saved selection to match the author's last edit location ---
$ws.Range("T3").Select() | Out-Null
